# Lecture 03 / 04 update:
#  1. Duplicate the "Today's objectives" slide (slide 5) and move the
#     duplicate to the end of the deck (new slide 8) before editing
#     anything, so the duplicate starts from the untouched original text.
#  2. On the new slide 8, insert two extra sub-bullets and drop the
#     trailing "Apply three different methods..." bullet.
#  3. On the original slide 5, just drop the trailing "Apply three
#     different methods..." bullet.

$p = $ppt.ActivePresentation

$sourceSlide = $p.Slides.Item(5)

# --- Step 1: duplicate slide 5 and move the copy to the end ---------------
$newSlide = $sourceSlide.Duplicate()
$newSlide.MoveTo($p.Slides.Count)

# --- Step 2: build out the new slide 8 content -----------------------------
$shp8 = $newSlide.Shapes.Item(2)
$tr8 = $shp8.TextFrame.TextRange

# Insert "String, int, float, bool," as a new sub-bullet right after the
# "Recite the main data types..." paragraph (paragraph 2).
$afterPara2 = $tr8.Paragraphs(2, 1)
$null = $afterPara2.InsertAfter("`rString, int, float, bool,")
$tr8.Paragraphs(3, 1).IndentLevel = 2

# Insert "List, dict, tuple, set" as a new sub-bullet right after the
# "Describe four different Python data structures..." paragraph (now
# paragraph 4, having shifted down by one from the previous insert).
$afterPara4 = $tr8.Paragraphs(4, 1)
$null = $afterPara4.InsertAfter("`rList, dict, tuple, set")
$tr8.Paragraphs(5, 1).IndentLevel = 2

# Drop the trailing "Apply three different methods..." paragraph while
# keeping the preceding "...dict of lists (Lecture 03.4)" paragraph intact.
$count8 = $tr8.Paragraphs().Count
$keepPara8 = $tr8.Paragraphs($count8 - 1, 1)
$dropPara8 = $tr8.Paragraphs($count8, 1)
$mergedLen8 = ($dropPara8.Start + $dropPara8.Length) - $keepPara8.Start
$mergedRange8 = $tr8.Characters($keepPara8.Start, $mergedLen8)
$mergedRange8.Text = "Describe Pandas DataFrames as a dict of lists (Lecture 03.4)"
$phantom8 = $tr8.Paragraphs($tr8.Paragraphs().Count, 1)
if ($phantom8.Length -eq 0) {
    $phantom8.Delete()
}

# Re-apply italics to "dict" and "lists " inside that last paragraph.
$lastPara8 = $tr8.Paragraphs($tr8.Paragraphs().Count, 1)
$lastText8 = $lastPara8.Text
$dictStart8 = $lastText8.IndexOf("dict") + 1
$lastPara8.Characters($dictStart8, 4).Font.Italic = $true
$listsStart8 = $lastText8.IndexOf("lists ") + 1
$lastPara8.Characters($listsStart8, 6).Font.Italic = $true

# --- Step 3: trim the trailing paragraph from the original slide 5 --------
$shp5 = $sourceSlide.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange
$count5 = $tr5.Paragraphs().Count
$keepPara5 = $tr5.Paragraphs($count5 - 1, 1)
$dropPara5 = $tr5.Paragraphs($count5, 1)
$mergedLen5 = ($dropPara5.Start + $dropPara5.Length) - $keepPara5.Start
$mergedRange5 = $tr5.Characters($keepPara5.Start, $mergedLen5)
$mergedRange5.Text = "Describe Pandas DataFrames as a dict of lists (Lecture 03.4)"
$phantom5 = $tr5.Paragraphs($tr5.Paragraphs().Count, 1)
if ($phantom5.Length -eq 0) {
    $phantom5.Delete()
}

$lastPara5 = $tr5.Paragraphs($tr5.Paragraphs().Count, 1)
$lastText5 = $lastPara5.Text
$dictStart5 = $lastText5.IndexOf("dict") + 1
$lastPara5.Characters($dictStart5, 4).Font.Italic = $true
$listsStart5 = $lastText5.IndexOf("lists ") + 1
$lastPara5.Characters($listsStart5, 6).Font.Italic = $true

Write-Host "Slides: $($p.Slides.Count)"
